$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (D1 stays "nama"; E1/F1/G1 renamed) ---
$ws.Range("D1").Value = "nama"
$ws.Range("E1").Value = "nilai_ujian_sekolah"
$ws.Range("F1").Value = "nilai_usp"
$ws.Range("G1").Value = "rerata"

# F1 previously carried a leftover date-number-format style; clear it back to Normal
$ws.Range("F1").Style = "Normal"

# --- Updated nilai_usp (column F) scores ---
$ws.Range("F2").Value = 100
$ws.Range("F3").Value = 86
$ws.Range("F4").Value = 79
$ws.Range("F5").Value = 77
$ws.Range("F6").Value = 85
$ws.Range("F7").Value = 92
$ws.Range("F8").Value = 78
$ws.Range("F9").Value = 90
$ws.Range("F10").Value = 88

# Column G ("rerata") already holds =AVERAGE(E:F) formulas, so they
# recalculate automatically once the F values above are written.

# --- Column widths: only column C keeps an explicit best-fit width now ---
$ws.Cells.ClearFormats()
$ws.Columns.Item(3).ColumnWidth = 14

# --- Selection moved to G3 ---
$ws.Range("G3").Select()
